$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed numeric error values (rows 2-15, columns B:K)
$ws.Range("B2").Value = 0.4657705279028737
$ws.Range("C2").Value = 0.07625703655992344
$ws.Range("D2").Value = 0.2445323695583424
$ws.Range("E2").Value = 0.7508699121205502
$ws.Range("F2").Value = 1.14353756626052
$ws.Range("G2").Value = 0.2224909929985871
$ws.Range("H2").Value = -0.5136315835481774
$ws.Range("I2").Value = 1.163476923566464
$ws.Range("J2").Value = 0.540088728946798
$ws.Range("K2").Value = 0.7747069251866952
$ws.Range("B3").Value = 0.09868501619293202
$ws.Range("C3").Value = 0.2471359091207911
$ws.Range("D3").Value = 0.728802226466535
$ws.Range("E3").Value = 1.137555585956165
$ws.Range("F3").Value = 0.217147742816176
$ws.Range("G3").Value = -0.5240001169185436
$ws.Range("H3").Value = 1.154999601442841
$ws.Range("I3").Value = 0.5323164206399872
$ws.Range("J3").Value = 0.7660829516537266
$ws.Range("K3").Value = 0.6652412546050546
$ws.Range("B4").Value = 0.2685604805780112
$ws.Range("C4").Value = 0.8489800752657474
$ws.Range("D4").Value = 1.048361241542382
$ws.Range("E4").Value = 0.1874207389675982
$ws.Range("F4").Value = -0.5153739968886335
$ws.Range("G4").Value = 1.129059428017406
$ws.Range("H4").Value = 0.5095378576117597
$ws.Range("I4").Value = 0.7524649165051087
$ws.Range("J4").Value = 0.6467204574198788
$ws.Range("K4").Value = -0.1454929044188731
$ws.Range("B5").Value = 0.806659442945358
$ws.Range("C5").Value = 1.019156767686649
$ws.Range("D5").Value = 0.184400139162677
$ws.Range("E5").Value = -0.534235280500114
$ws.Range("F5").Value = 1.108728296952201
$ws.Range("G5").Value = 0.4950077021775119
$ws.Range("H5").Value = 0.7358238710128409
$ws.Range("I5").Value = 0.6291366119911286
$ws.Range("J5").Value = -0.1619667879489159
$ws.Range("K5").Value = 0.449372724506711
$ws.Range("B6").Value = 1.359354508304559
$ws.Range("C6").Value = 0.2597712009466141
$ws.Range("D6").Value = -0.7253285668131131
$ws.Range("E6").Value = 1.133422209538934
$ws.Range("F6").Value = 0.5043871057338079
$ws.Range("G6").Value = 0.6810130724014498
$ws.Range("H6").Value = 0.6102425211024588
$ws.Range("I6").Value = -0.1768307887639616
$ws.Range("J6").Value = 0.4206610138773402
$ws.Range("K6").Value = 0.2638965897873631
$ws.Range("B7").Value = 0.7104660729368646
$ws.Range("C7").Value = -0.6782178646007546
$ws.Range("D7").Value = 0.8931389787481262
$ws.Range("E7").Value = 0.5403321867996707
$ws.Range("F7").Value = 0.6789623149618569
$ws.Range("G7").Value = 0.5336742816664286
$ws.Range("H7").Value = -0.204230183139245
$ws.Range("I7").Value = 0.3943719636796149
$ws.Range("J7").Value = 0.2204992990740305
$ws.Range("B8").Value = -0.3658922776772162
$ws.Range("C8").Value = 1.026427718483651
$ws.Range("D8").Value = 0.3612537649834791
$ws.Range("E8").Value = 0.707156553732424
$ws.Range("F8").Value = 0.5696995928588221
$ws.Range("G8").Value = -0.2417717533434586
$ws.Range("H8").Value = 0.3888296245922537
$ws.Range("I8").Value = 0.2246746280127792
$ws.Range("B9").Value = 1.262018209591492
$ws.Range("C9").Value = 0.4459325058577887
$ws.Range("D9").Value = 0.5615561870987069
$ws.Range("E9").Value = 0.5796533357180647
$ws.Range("F9").Value = -0.2264065836439137
$ws.Range("G9").Value = 0.3499812896348306
$ws.Range("H9").Value = 0.2088288189855932
$ws.Range("B10").Value = 0.7569566923391715
$ws.Range("C10").Value = 0.6786449615099022
$ws.Range("D10").Value = 0.4171891942684979
$ws.Range("E10").Value = -0.1974476331787121
$ws.Range("F10").Value = 0.3857241620897341
$ws.Range("G10").Value = 0.1775011726019661
$ws.Range("B11").Value = 0.9254701389140165
$ws.Range("C11").Value = 0.4348450618063874
$ws.Range("D11").Value = -0.2920698722897066
$ws.Range("E11").Value = 0.4179415503382142
$ws.Range("F11").Value = 0.1898892984296834
$ws.Range("B12").Value = 0.6745214212225993
$ws.Range("C12").Value = -0.2070456288204931
$ws.Range("D12").Value = 0.3014569719802002
$ws.Range("E12").Value = 0.2049945700815359
$ws.Range("B13").Value = -0.04218555178640582
$ws.Range("C13").Value = 0.3149942442281164
$ws.Range("D13").Value = 0.1420216510915729
$ws.Range("B14").Value = 0.5688432860935244
$ws.Range("C14").Value = 0.2413397012736094
$ws.Range("B15").Value = 0.2853993925130583

# Remove cells that are no longer part of the shrinking staircase range
$ws.Range("K7").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("H10").ClearContents()
$ws.Range("G11").ClearContents()
$ws.Range("F12").ClearContents()
$ws.Range("E13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("B16").ClearContents()
